# Applies:
#  1. Inserts a new "wiring" slide (Title + Content layout) before the
#     "software" slide -> new slide ends up at position 5, pushing the
#     existing "software"/"interface" slides down one spot.
#  2. Renames the "interface" slide's title to "What's next".

$p = $ppt.ActivePresentation

# --- 1. Insert the new "wiring" slide at position 5 -------------------
# Slide 5 currently holds "software" (Title and Content layout); reuse
# that same layout for the new slide so shapes/placeholders match.
$layout = $p.Slides.Item(5).CustomLayout
$newSlide = $p.Slides.AddSlide(5, $layout)
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "wiring"

# --- 2. Rename the "interface" slide's title to "What's next" --------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $sl = $p.Slides.Item($i)
    if ($sl.Shapes.Item(1).TextFrame.TextRange.Text -eq "interface") {
        $sl.Shapes.Item(1).TextFrame.TextRange.Text = "What" + [char]0x2019 + "s next"
        break
    }
}
